$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value. Numeric-looking "Price" strings are
# prefixed with a leading apostrophe so Excel stores them as text (matching
# the source data's display-formatted strings) instead of auto-converting
# them to numbers.
$updates = @(
    @{ Cell = 'D2'; Value = '30.787.62' }
    @{ Cell = 'D3'; Value = '1.937.18' }
    @{ Cell = 'E3'; Value = '  -0.75%  ' }
    @{ Cell = 'D4'; Value = "'1.000" }
    @{ Cell = 'E4'; Value = '  +0.04%  ' }
    @{ Cell = 'D5'; Value = "'243.42" }
    @{ Cell = 'E5'; Value = '  -0.79%  ' }
    @{ Cell = 'D6'; Value = "'1.000" }
    @{ Cell = 'E6'; Value = '  +0.04%  ' }
    @{ Cell = 'D7'; Value = "'0.4880" }
    @{ Cell = 'E7'; Value = '  -0.05%  ' }
    @{ Cell = 'D8'; Value = "'0.2953" }
    @{ Cell = 'E8'; Value = '  -0.35%  ' }
    @{ Cell = 'D9'; Value = "'0.06894" }
    @{ Cell = 'E9'; Value = '  +0.96%  ' }
    @{ Cell = 'D10'; Value = "'19.34" }
    @{ Cell = 'E10'; Value = '  +0.96%  ' }
    @{ Cell = 'D11'; Value = "'105.04" }
    @{ Cell = 'E11'; Value = '  -1.93%  ' }
    @{ Cell = 'D12'; Value = '1.942.17' }
    @{ Cell = 'E12'; Value = '  -0.70%  ' }
    @{ Cell = 'D13'; Value = "'0.07781" }
    @{ Cell = 'E13'; Value = '  +0.52%  ' }
    @{ Cell = 'D14'; Value = "'5.352" }
    @{ Cell = 'E14'; Value = '  -2.11%  ' }
    @{ Cell = 'D15'; Value = "'0.7025" }
    @{ Cell = 'E15'; Value = '  -0.53%  ' }
    @{ Cell = 'D16'; Value = "'273.19" }
    @{ Cell = 'E16'; Value = '  -3.00%  ' }
    @{ Cell = 'D17'; Value = '30.800.39' }
    @{ Cell = 'E17'; Value = '  -0.66%  ' }
    @{ Cell = 'D18'; Value = "'5.704" }
    @{ Cell = 'E18'; Value = '  +3.69%  ' }
    @{ Cell = 'D19'; Value = "'0.000007731" }
    @{ Cell = 'E19'; Value = '  +0.01%  ' }
    @{ Cell = 'D20'; Value = "'13.11" }
    @{ Cell = 'E20'; Value = '  -0.89%  ' }
    @{ Cell = 'D21'; Value = "'0.9997" }
    @{ Cell = 'E21'; Value = '  +0.05%  ' }
    @{ Cell = 'D22'; Value = "'0.9999" }
    @{ Cell = 'E22'; Value = '  +0.09%  ' }
    @{ Cell = 'D23'; Value = "'6.538" }
    @{ Cell = 'E23'; Value = '  +0.66%  ' }
    @{ Cell = 'D24'; Value = "'9.812" }
    @{ Cell = 'E24'; Value = '  -0.30%  ' }
    @{ Cell = 'D25'; Value = "'164.97" }
    @{ Cell = 'E25'; Value = '  -2.57%  ' }
    @{ Cell = 'E26'; Value = '  -1.84%  ' }
    @{ Cell = 'D27'; Value = "'2.163" }
    @{ Cell = 'E27'; Value = '  -2.24%  ' }
    @{ Cell = 'D28'; Value = "'0.1036" }
    @{ Cell = 'E28'; Value = '  -1.73%  ' }
    @{ Cell = 'D29'; Value = "'1.386" }
    @{ Cell = 'E29'; Value = '  -2.57%  ' }
    @{ Cell = 'D30'; Value = "'4.694" }
    @{ Cell = 'E30'; Value = '  +2.75%  ' }
    @{ Cell = 'E31'; Value = '  -1.34%  ' }
    @{ Cell = 'D32'; Value = "'4.430" }
    @{ Cell = 'E32'; Value = '  -0.75%  ' }
    @{ Cell = 'D33'; Value = "'0.04912" }
    @{ Cell = 'E33'; Value = '  -0.80%  ' }
    @{ Cell = 'D34'; Value = "'0.7608" }
    @{ Cell = 'E34'; Value = '  -1.18%  ' }
    @{ Cell = 'D35'; Value = "'1.152" }
    @{ Cell = 'E35'; Value = '  -1.48%  ' }
    @{ Cell = 'D36'; Value = "'0.9992" }
    @{ Cell = 'E36'; Value = '  +0.05%  ' }
    @{ Cell = 'E37'; Value = '  +0.09%  ' }
    @{ Cell = 'D38'; Value = "'0.02010" }
    @{ Cell = 'E38'; Value = '  -0.45%  ' }
    @{ Cell = 'D39'; Value = "'79.31" }
    @{ Cell = 'E39'; Value = '  +6.38%  ' }
    @{ Cell = 'D40'; Value = "'2.667" }
    @{ Cell = 'E40'; Value = '  -1.17%  ' }
    @{ Cell = 'D41'; Value = "'6.493" }
    @{ Cell = 'E41'; Value = '  -0.64%  ' }
    @{ Cell = 'D42'; Value = "'2.083" }
    @{ Cell = 'E42'; Value = '  -3.41%  ' }
    @{ Cell = 'D43'; Value = "'0.9016" }
    @{ Cell = 'E43'; Value = '  +2.08%  ' }
    @{ Cell = 'D44'; Value = "'0.4456" }
    @{ Cell = 'E44'; Value = '  -0.82%  ' }
    @{ Cell = 'E45'; Value = '  -0.80%  ' }
    @{ Cell = 'D46'; Value = "'7.892" }
    @{ Cell = 'E46'; Value = '  -3.06%  ' }
    @{ Cell = 'D47'; Value = "'1.000" }
    @{ Cell = 'E47'; Value = '  +0.07%  ' }
    @{ Cell = 'D48'; Value = "'992.33" }
    @{ Cell = 'E48'; Value = '  +1.26%  ' }
    @{ Cell = 'E49'; Value = '  -1.03%  ' }
    @{ Cell = 'D50'; Value = "'36.31" }
    @{ Cell = 'E50'; Value = '  +1.43%  ' }
    @{ Cell = 'D51'; Value = "'9.231" }
    @{ Cell = 'E51'; Value = '  -1.72%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
